# Auto-generated Excel COM-interop script to apply the XML diff changes
# to Sheets/Raiden_Profits.xlsx (workbook with sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 47654  # H3: 46769.332 -> 47654
$ws.Cells.Item(3, 10).Value = 47654  # J3: 46769.332 -> 47654
$ws.Cells.Item(3, 12).Value = 47654  # L3: 46769.332 -> 47654
$ws.Cells.Item(3, 14).Value = -47882  # N3: -46997.332 -> -47882
$ws.Cells.Item(52, 8).Value = 350  # H52: 0 -> 350
$ws.Cells.Item(52, 9).Value = 350  # I52: 0 -> 350
$ws.Cells.Item(52, 11).Value = 1050  # K52: 0 -> 1050
$ws.Cells.Item(52, 13).Value = -890  # M52: __ABSENT__ -> -890
$ws.Cells.Item(53, 8).Value = 150.27272  # H53: 165.2 -> 150.27272
$ws.Cells.Item(53, 9).Value = 53  # I53: 63.4 -> 53
$ws.Cells.Item(53, 11).Value = 53  # K53: 63.4 -> 53
$ws.Cells.Item(53, 13).Value = 584  # M53: 573.6 -> 584
$ws.Cells.Item(80, 8).Value = 862.75  # H80: 871.7143 -> 862.75
$ws.Cells.Item(80, 9).Value = 649.5  # I80: 733 -> 649.5
$ws.Cells.Item(80, 10).Value = 933.8333  # J80: 975.75 -> 933.8333
$ws.Cells.Item(80, 11).Value = 1948.5  # K80: 2199 -> 1948.5
$ws.Cells.Item(80, 12).Value = 2801.4999  # L80: 2927.25 -> 2801.4999
$ws.Cells.Item(80, 13).Value = -950.5  # M80: -1201 -> -950.5
$ws.Cells.Item(80, 14).Value = -4797.4999  # N80: -4923.25 -> -4797.4999
$ws.Cells.Item(83, 8).Value = 862.75  # H83: 871.7143 -> 862.75
$ws.Cells.Item(83, 9).Value = 649.5  # I83: 733 -> 649.5
$ws.Cells.Item(83, 10).Value = 933.8333  # J83: 975.75 -> 933.8333
$ws.Cells.Item(83, 11).Value = 5845.5  # K83: 6597 -> 5845.5
$ws.Cells.Item(83, 12).Value = 8404.4997  # L83: 8781.75 -> 8404.4997
$ws.Cells.Item(83, 13).Value = -853.5  # M83: -1605 -> -853.5
$ws.Cells.Item(83, 14).Value = -18388.4997  # N83: -18765.75 -> -18388.4997
$ws.Cells.Item(98, 8).Value = 1072.875  # H98: 1072.9375 -> 1072.875
$ws.Cells.Item(98, 9).Value = 1072.875  # I98: 1095.6 -> 1072.875
$ws.Cells.Item(98, 10).Value = 0  # J98: 733 -> 0
$ws.Cells.Item(98, 11).Value = 1072.875  # K98: 1095.6 -> 1072.875
$ws.Cells.Item(98, 12).Value = 0  # L98: 733 -> 0
$ws.Cells.Item(98, 13).Value = 425.125  # M98: 402.4000000000001 -> 425.125
$ws.Cells.Item(98, 14).ClearContents()  # N98: -3729 -> (removed)
$ws.Cells.Item(102, 8).Value = 47654  # H102: 46769.332 -> 47654
$ws.Cells.Item(102, 10).Value = 47654  # J102: 46769.332 -> 47654
$ws.Cells.Item(102, 12).Value = 47654  # L102: 46769.332 -> 47654
$ws.Cells.Item(102, 14).Value = -54144  # N102: -53259.332 -> -54144
$ws.Cells.Item(112, 8).Value = 3137.4333  # H112: 3174.3794 -> 3137.4333
$ws.Cells.Item(112, 9).Value = 1300  # I112: 0 -> 1300
$ws.Cells.Item(112, 10).Value = 3200.7932  # J112: 3174.3794 -> 3200.7932
$ws.Cells.Item(112, 11).Value = 3900  # K112: 0 -> 3900
$ws.Cells.Item(112, 12).Value = 9602.3796  # L112: 9523.138199999999 -> 9602.3796
$ws.Cells.Item(112, 13).Value = -2792  # M112: __ABSENT__ -> -2792
$ws.Cells.Item(112, 14).Value = -11818.3796  # N112: -11739.1382 -> -11818.3796
$ws.Cells.Item(122, 8).Value = 1072.875  # H122: 1072.9375 -> 1072.875
$ws.Cells.Item(122, 9).Value = 1072.875  # I122: 1095.6 -> 1072.875
$ws.Cells.Item(122, 10).Value = 0  # J122: 733 -> 0
$ws.Cells.Item(122, 11).Value = 3218.625  # K122: 3286.8 -> 3218.625
$ws.Cells.Item(122, 12).Value = 0  # L122: 2199 -> 0
$ws.Cells.Item(122, 13).Value = -768.625  # M122: -836.7999999999997 -> -768.625
$ws.Cells.Item(122, 14).ClearContents()  # N122: -7099 -> (removed)
$ws.Cells.Item(138, 8).Value = 2661.5  # H138: 2683.7942 -> 2661.5
$ws.Cells.Item(138, 9).Value = 2003.6666  # I138: 1957.7727 -> 2003.6666
$ws.Cells.Item(138, 10).Value = 2943.4285  # J138: 3031.0217 -> 2943.4285
$ws.Cells.Item(138, 11).Value = 6010.9998  # K138: 5873.3181 -> 6010.9998
$ws.Cells.Item(138, 12).Value = 8830.2855  # L138: 9093.0651 -> 8830.2855
$ws.Cells.Item(138, 13).Value = -870.9997999999996  # M138: -733.3181000000004 -> -870.9997999999996
$ws.Cells.Item(138, 14).Value = -19110.2855  # N138: -19373.0651 -> -19110.2855
$ws.Cells.Item(141, 8).Value = 4705.1763  # H141: 4217.9473 -> 4705.1763
$ws.Cells.Item(141, 9).Value = 4430.3335  # I141: 3613.5454 -> 4430.3335
$ws.Cells.Item(141, 10).Value = 5014.375  # J141: 5049 -> 5014.375
$ws.Cells.Item(141, 11).Value = 13291.0005  # K141: 10840.6362 -> 13291.0005
$ws.Cells.Item(141, 12).Value = 15043.125  # L141: 15147 -> 15043.125
$ws.Cells.Item(141, 13).Value = -8111.000499999998  # M141: -5660.636200000001 -> -8111.000499999998
$ws.Cells.Item(141, 14).Value = -25403.125  # N141: -25507 -> -25403.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17748.887  # H32: 18269.77 -> 17748.887
$ws.Cells.Item(32, 9).Value = 6013.804  # I32: 6459.1113 -> 6013.804
$ws.Cells.Item(32, 11).Value = 6013.804  # K32: 6459.1113 -> 6013.804
$ws.Cells.Item(32, 13).Value = -5726.804  # M32: -6172.1113 -> -5726.804
$ws.Cells.Item(63, 8).Value = 2722.0476  # H63: 2692.465 -> 2722.0476
$ws.Cells.Item(63, 9).Value = 2697.3333  # I63: 2663.6216 -> 2697.3333
$ws.Cells.Item(63, 11).Value = 2697.3333  # K63: 2663.6216 -> 2697.3333
$ws.Cells.Item(63, 13).Value = -2011.3333  # M63: -1977.6216 -> -2011.3333
$ws.Cells.Item(66, 8).Value = 2722.0476  # H66: 2692.465 -> 2722.0476
$ws.Cells.Item(66, 9).Value = 2697.3333  # I66: 2663.6216 -> 2697.3333
$ws.Cells.Item(66, 11).Value = 13486.6665  # K66: 13318.108 -> 13486.6665
$ws.Cells.Item(66, 13).Value = -10054.6665  # M66: -9886.108 -> -10054.6665
$ws.Cells.Item(122, 8).Value = 3041.054  # H122: 3053.5715 -> 3041.054
$ws.Cells.Item(122, 9).Value = 2970.6667  # I122: 2961.2068 -> 2970.6667
$ws.Cells.Item(122, 10).Value = 3342.7144  # J122: 3500 -> 3342.7144
$ws.Cells.Item(122, 11).Value = 8912.000100000001  # K122: 8883.6204 -> 8912.000100000001
$ws.Cells.Item(122, 12).Value = 10028.1432  # L122: 10500 -> 10028.1432
$ws.Cells.Item(122, 13).Value = -6462.000100000001  # M122: -6433.6204 -> -6462.000100000001
$ws.Cells.Item(122, 14).Value = -14928.1432  # N122: -15400 -> -14928.1432
$ws.Cells.Item(132, 8).Value = 3158.5557  # H132: 2380.1428 -> 3158.5557
$ws.Cells.Item(132, 9).Value = 1204.5714  # I132: 1110.5834 -> 1204.5714
$ws.Cells.Item(132, 11).Value = 3613.7142  # K132: 3331.7502 -> 3613.7142
$ws.Cells.Item(132, 13).Value = -1083.7142  # M132: -801.7501999999999 -> -1083.7142

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2339.1875  # H20: 2412.8572 -> 2339.1875
$ws.Cells.Item(20, 9).Value = 2348.5386  # I20: 2402.6667 -> 2348.5386
$ws.Cells.Item(20, 10).Value = 2298.6667  # J20: 2474 -> 2298.6667
$ws.Cells.Item(20, 11).Value = 2348.5386  # K20: 2402.6667 -> 2348.5386
$ws.Cells.Item(20, 12).Value = 2298.6667  # L20: 2474 -> 2298.6667
$ws.Cells.Item(20, 13).Value = -2101.5386  # M20: -2155.6667 -> -2101.5386
$ws.Cells.Item(20, 14).Value = -2792.6667  # N20: -2968 -> -2792.6667
$ws.Cells.Item(99, 8).Value = 1297.9333  # H99: 1415.3846 -> 1297.9333
$ws.Cells.Item(99, 9).Value = 1298.5555  # I99: 1451 -> 1298.5555
$ws.Cells.Item(99, 10).Value = 1297  # J99: 1358.4 -> 1297
$ws.Cells.Item(99, 11).Value = 1298.5555  # K99: 1451 -> 1298.5555
$ws.Cells.Item(99, 12).Value = 1297  # L99: 1358.4 -> 1297
$ws.Cells.Item(99, 13).Value = 199.4445000000001  # M99: 47 -> 199.4445000000001
$ws.Cells.Item(99, 14).Value = -4293  # N99: -4354.4 -> -4293
$ws.Cells.Item(105, 8).Value = 3539.842  # H105: 3631 -> 3539.842
$ws.Cells.Item(105, 9).Value = 2988.875  # I105: 3144.5715 -> 2988.875
$ws.Cells.Item(105, 11).Value = 2988.875  # K105: 3144.5715 -> 2988.875
$ws.Cells.Item(105, 13).Value = -1241.875  # M105: -1397.5715 -> -1241.875
$ws.Cells.Item(107, 8).Value = 1871.3334  # H107: 1927.5 -> 1871.3334
$ws.Cells.Item(107, 9).Value = 1905.6  # I107: 1927.5 -> 1905.6
$ws.Cells.Item(107, 10).Value = 1700  # J107: 0 -> 1700
$ws.Cells.Item(107, 11).Value = 1905.6  # K107: 1927.5 -> 1905.6
$ws.Cells.Item(107, 12).Value = 1700  # L107: 0 -> 1700
$ws.Cells.Item(107, 13).Value = 14.40000000000009  # M107: -7.5 -> 14.40000000000009
$ws.Cells.Item(107, 14).Value = -5540  # N107: __ABSENT__ -> -5540

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1509.4  # H16: 1665.2222 -> 1509.4
$ws.Cells.Item(16, 10).Value = 1227.2  # J16: 1507.25 -> 1227.2
$ws.Cells.Item(16, 12).Value = 1227.2  # L16: 1507.25 -> 1227.2
$ws.Cells.Item(16, 14).Value = -1801.2  # N16: -2081.25 -> -1801.2
$ws.Cells.Item(31, 8).Value = 5941.7544  # H31: 6360.8335 -> 5941.7544
$ws.Cells.Item(31, 9).Value = 3246.842  # I31: 3512.5 -> 3246.842
$ws.Cells.Item(31, 10).Value = 11331.579  # J31: 10633.333 -> 11331.579
$ws.Cells.Item(31, 11).Value = 3246.842  # K31: 3512.5 -> 3246.842
$ws.Cells.Item(31, 12).Value = 11331.579  # L31: 10633.333 -> 11331.579
$ws.Cells.Item(31, 13).Value = -2951.842  # M31: -3217.5 -> -2951.842
$ws.Cells.Item(31, 14).Value = -11921.579  # N31: -11223.333 -> -11921.579
$ws.Cells.Item(34, 8).Value = 5941.7544  # H34: 6360.8335 -> 5941.7544
$ws.Cells.Item(34, 9).Value = 3246.842  # I34: 3512.5 -> 3246.842
$ws.Cells.Item(34, 10).Value = 11331.579  # J34: 10633.333 -> 11331.579
$ws.Cells.Item(34, 11).Value = 3246.842  # K34: 3512.5 -> 3246.842
$ws.Cells.Item(34, 12).Value = 11331.579  # L34: 10633.333 -> 11331.579
$ws.Cells.Item(34, 13).Value = -3044.842  # M34: -3310.5 -> -3044.842
$ws.Cells.Item(34, 14).Value = -11735.579  # N34: -11037.333 -> -11735.579
$ws.Cells.Item(58, 8).Value = 2138.52  # H58: 2372.6667 -> 2138.52
$ws.Cells.Item(58, 9).Value = 2243.15  # I58: 2483.2354 -> 2243.15
$ws.Cells.Item(58, 10).Value = 1720  # J58: 1902.75 -> 1720
$ws.Cells.Item(58, 11).Value = 2243.15  # K58: 2483.2354 -> 2243.15
$ws.Cells.Item(58, 12).Value = 1720  # L58: 1902.75 -> 1720
$ws.Cells.Item(58, 13).Value = -2040.15  # M58: -2280.2354 -> -2040.15
$ws.Cells.Item(58, 14).Value = -2126  # N58: -2308.75 -> -2126
$ws.Cells.Item(99, 9).Value = 1200  # I99: 1199.5 -> 1200
$ws.Cells.Item(99, 10).Value = 1199  # J99: 0 -> 1199
$ws.Cells.Item(99, 11).Value = 1200  # K99: 1199.5 -> 1200
$ws.Cells.Item(99, 12).Value = 1199  # L99: 0 -> 1199
$ws.Cells.Item(99, 13).Value = 298  # M99: 298.5 -> 298
$ws.Cells.Item(99, 14).Value = -4195  # N99: __ABSENT__ -> -4195
$ws.Cells.Item(113, 8).Value = 1509.4  # H113: 1665.2222 -> 1509.4
$ws.Cells.Item(113, 10).Value = 1227.2  # J113: 1507.25 -> 1227.2
$ws.Cells.Item(113, 12).Value = 1227.2  # L113: 1507.25 -> 1227.2
$ws.Cells.Item(113, 14).Value = -5567.2  # N113: -5847.25 -> -5567.2
$ws.Cells.Item(126, 9).Value = 1200  # I126: 1199.5 -> 1200
$ws.Cells.Item(126, 10).Value = 1199  # J126: 0 -> 1199
$ws.Cells.Item(126, 11).Value = 3600  # K126: 3598.5 -> 3600
$ws.Cells.Item(126, 12).Value = 3597  # L126: 0 -> 3597
$ws.Cells.Item(126, 13).Value = -1130  # M126: -1128.5 -> -1130
$ws.Cells.Item(126, 14).Value = -8537  # N126: __ABSENT__ -> -8537
$ws.Cells.Item(136, 8).Value = 2138.52  # H136: 2372.6667 -> 2138.52
$ws.Cells.Item(136, 9).Value = 2243.15  # I136: 2483.2354 -> 2243.15
$ws.Cells.Item(136, 10).Value = 1720  # J136: 1902.75 -> 1720
$ws.Cells.Item(136, 11).Value = 6729.450000000001  # K136: 7449.706200000001 -> 6729.450000000001
$ws.Cells.Item(136, 12).Value = 5160  # L136: 5708.25 -> 5160
$ws.Cells.Item(136, 13).Value = -4179.450000000001  # M136: -4899.706200000001 -> -4179.450000000001
$ws.Cells.Item(136, 14).Value = -10260  # N136: -10808.25 -> -10260

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 112.48148  # H2: 117.12 -> 112.48148
$ws.Cells.Item(2, 9).Value = 128.09091  # I2: 135.45 -> 128.09091
$ws.Cells.Item(2, 11).Value = 128.09091  # K2: 135.45 -> 128.09091
$ws.Cells.Item(2, 13).Value = -15.09091000000001  # M2: -22.44999999999999 -> -15.09091000000001
$ws.Cells.Item(80, 8).Value = 7718.75  # H80: 9115.462 -> 7718.75
$ws.Cells.Item(80, 9).Value = 2418.4546  # I80: 2700.5 -> 2418.4546
$ws.Cells.Item(80, 11).Value = 2418.4546  # K80: 2700.5 -> 2418.4546
$ws.Cells.Item(80, 13).Value = -1420.4546  # M80: -1702.5 -> -1420.4546
$ws.Cells.Item(83, 8).Value = 7718.75  # H83: 9115.462 -> 7718.75
$ws.Cells.Item(83, 9).Value = 2418.4546  # I83: 2700.5 -> 2418.4546
$ws.Cells.Item(83, 11).Value = 12092.273  # K83: 13502.5 -> 12092.273
$ws.Cells.Item(83, 13).Value = -7100.273000000001  # M83: -8510.5 -> -7100.273000000001
$ws.Cells.Item(132, 8).Value = 3291.0588  # H132: 2369.111 -> 3291.0588
$ws.Cells.Item(132, 9).Value = 4826.5  # I132: 2311.0625 -> 4826.5
$ws.Cells.Item(132, 11).Value = 14479.5  # K132: 6933.1875 -> 14479.5
$ws.Cells.Item(132, 13).Value = -11949.5  # M132: -4403.1875 -> -11949.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5366.7896  # H40: 5348.3 -> 5366.7896
$ws.Cells.Item(40, 9).Value = 5999  # I40: 5873.75 -> 5999
$ws.Cells.Item(40, 11).Value = 5999  # K40: 5873.75 -> 5999
$ws.Cells.Item(40, 13).Value = -5863  # M40: -5737.75 -> -5863
$ws.Cells.Item(55, 8).Value = 535.2222  # H55: 555.6667 -> 535.2222
$ws.Cells.Item(55, 9).Value = 736.8333  # I55: 767.5 -> 736.8333
$ws.Cells.Item(55, 11).Value = 736.8333  # K55: 767.5 -> 736.8333
$ws.Cells.Item(55, 13).Value = -563.8333  # M55: -594.5 -> -563.8333
$ws.Cells.Item(122, 8).Value = 3615.1538  # H122: 3615.2307 -> 3615.1538
$ws.Cells.Item(122, 9).Value = 2698  # I122: 0 -> 2698
$ws.Cells.Item(122, 10).Value = 3691.5833  # J122: 3615.2307 -> 3691.5833
$ws.Cells.Item(122, 11).Value = 8094  # K122: 0 -> 8094
$ws.Cells.Item(122, 12).Value = 11074.7499  # L122: 10845.6921 -> 11074.7499
$ws.Cells.Item(122, 13).Value = -5644  # M122: __ABSENT__ -> -5644
$ws.Cells.Item(122, 14).Value = -15974.7499  # N122: -15745.6921 -> -15974.7499

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 13302.24  # H122: 13823.083 -> 13302.24
$ws.Cells.Item(122, 9).Value = 11893.826  # I122: 12398 -> 11893.826
$ws.Cells.Item(122, 11).Value = 35681.478  # K122: 37194 -> 35681.478
$ws.Cells.Item(122, 13).Value = -33231.478  # M122: -34744 -> -33231.478
$ws.Cells.Item(136, 8).Value = 2764.5  # H136: 2805.1052 -> 2764.5
$ws.Cells.Item(136, 9).Value = 2771.7222  # I136: 2817.5293 -> 2771.7222
$ws.Cells.Item(136, 11).Value = 8315.1666  # K136: 8452.5879 -> 8315.1666
$ws.Cells.Item(136, 13).Value = -5765.1666  # M136: -5902.5879 -> -5765.1666
